$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing "Сдали комплект" (E column) marks for rows 9, 18 and 23.
$ws.Range("E9").Value = 5
$ws.Range("E18").Value = 5

# Row 23 previously had an empty D cell (ПП) - student now also passed it,
# plus gets the "Сдали комплект" mark in E.
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 5

# Move the active selection from E16 to E10 (no data change, just the
# last-selected cell as saved in the sheet view).
$ws.Range("E10").Select()
